$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 57379.05
$ws.Range("J32").Value = 69636.75
$ws.Range("L32").Value = 69636.75
$ws.Range("N32").Value = -70288.75

$ws.Range("H40").Value = 13623.857
$ws.Range("I40").Value = 6900
$ws.Range("J40").Value = 16313.4
$ws.Range("K40").Value = 6900
$ws.Range("L40").Value = 16313.4
$ws.Range("M40").Value = -6725
$ws.Range("N40").Value = -16663.4

$ws.Range("H80").Value = 2829.9614
$ws.Range("I80").Value = 1837.7646
$ws.Range("J80").Value = 4704.1113
$ws.Range("K80").Value = 5513.293799999999
$ws.Range("L80").Value = 14112.3339
$ws.Range("M80").Value = -4515.293799999999
$ws.Range("N80").Value = -16108.3339

$ws.Range("H83").Value = 2829.9614
$ws.Range("I83").Value = 1837.7646
$ws.Range("J83").Value = 4704.1113
$ws.Range("K83").Value = 16539.8814
$ws.Range("L83").Value = 42337.00169999999
$ws.Range("M83").Value = -11547.8814
$ws.Range("N83").Value = -52321.00169999999

$ws.Range("H86").Value = 63495536
$ws.Range("I86").Value = 30304486
$ws.Range("K86").Value = 30304486
$ws.Range("M86").Value = -30303363

$ws.Range("H89").Value = 63495536
$ws.Range("I89").Value = 30304486
$ws.Range("K89").Value = 151522430
$ws.Range("M89").Value = -151516814

$ws.Range("H132").Value = 1835.1086
$ws.Range("I132").Value = 1515.8948
$ws.Range("K132").Value = 4547.6844
$ws.Range("M132").Value = -2017.6844

$ws.Range("H137").Value = 485348.56
$ws.Range("I137").Value = 2197.4167
$ws.Range("J137").Value = 2417953.2
$ws.Range("K137").Value = 6592.250100000001
$ws.Range("L137").Value = 7253859.600000001
$ws.Range("M137").Value = -4042.250100000001
$ws.Range("N137").Value = -7258959.600000001

$ws.Range("H138").Value = 2758.9714
$ws.Range("I138").Value = 1446
$ws.Range("J138").Value = 3999
$ws.Range("K138").Value = 4338
$ws.Range("L138").Value = 11997
$ws.Range("M138").Value = 802
$ws.Range("N138").Value = -22277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11131.788
$ws.Range("I32").Value = 5798.5415
$ws.Range("K32").Value = 5798.5415
$ws.Range("M32").Value = -5511.5415

$ws.Range("H114").Value = 73000
$ws.Range("J114").Value = 73000
$ws.Range("L114").Value = 73000
$ws.Range("N114").Value = -81678

$ws.Range("H132").Value = 2358.8965
$ws.Range("I132").Value = 1767.0555
$ws.Range("J132").Value = 3327.3635
$ws.Range("K132").Value = 5301.166499999999
$ws.Range("L132").Value = 9982.0905
$ws.Range("M132").Value = -2771.166499999999
$ws.Range("N132").Value = -15042.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1873.8889
$ws.Range("I20").Value = 1222.875
$ws.Range("J20").Value = 2394.7
$ws.Range("K20").Value = 1222.875
$ws.Range("L20").Value = 2394.7
$ws.Range("M20").Value = -975.875
$ws.Range("N20").Value = -2888.7

$ws.Range("H86").Value = 1885.875
$ws.Range("I86").Value = 1815.909
$ws.Range("K86").Value = 1815.909
$ws.Range("M86").Value = -692.9090000000001

$ws.Range("H89").Value = 1885.875
$ws.Range("I89").Value = 1815.909
$ws.Range("K89").Value = 9079.545
$ws.Range("M89").Value = -3463.545

$ws.Range("H94").Value = 3123.2
$ws.Range("I94").Value = 3123.2
$ws.Range("K94").Value = 3123.2
$ws.Range("M94").Value = -2672.2

$ws.Range("H105").Value = 102853.8
$ws.Range("I105").Value = 334047.34
$ws.Range("J105").Value = 3770.8572
$ws.Range("K105").Value = 334047.34
$ws.Range("L105").Value = 3770.8572
$ws.Range("M105").Value = -332300.34
$ws.Range("N105").Value = -7264.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2058.3713
$ws.Range("I31").Value = 1504.9
$ws.Range("K31").Value = 1504.9
$ws.Range("M31").Value = -1209.9

$ws.Range("H34").Value = 2058.3713
$ws.Range("I34").Value = 1504.9
$ws.Range("K34").Value = 1504.9
$ws.Range("M34").Value = -1302.9

$ws.Range("H60").Value = 27214.143
$ws.Range("J60").Value = 37500
$ws.Range("L60").Value = 37500
$ws.Range("N60").Value = -38522

$ws.Range("H99").Value = 2607804
$ws.Range("I99").Value = 3599.8572
$ws.Range("K99").Value = 3599.8572
$ws.Range("M99").Value = -2101.8572

$ws.Range("H126").Value = 2607804
$ws.Range("I126").Value = 3599.8572
$ws.Range("K126").Value = 10799.5716
$ws.Range("M126").Value = -8329.571599999999

$ws.Range("H132").Value = 2055.2354
$ws.Range("I132").Value = 1593.9
$ws.Range("K132").Value = 4781.700000000001
$ws.Range("M132").Value = -2251.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 5151.25
$ws.Range("I7").Value = 34.666668
$ws.Range("K7").Value = 104.000004
$ws.Range("M7").Value = 7.999995999999996

$ws.Range("H8").Value = 346666.34
$ws.Range("I8").Value = 346666.34
$ws.Range("K8").Value = 1039999.02
$ws.Range("M8").Value = -1039860.02

$ws.Range("H104").Value = 4671.65
$ws.Range("I104").Value = 2978
$ws.Range("K104").Value = 8934
$ws.Range("M104").Value = -6313

$ws.Range("H107").Value = 609.25
$ws.Range("I107").Value = 663.5
$ws.Range("K107").Value = 1990.5
$ws.Range("M107").Value = -70.5

$ws.Range("H116").Value = 1258.1428
$ws.Range("I116").Value = 1231.3334
$ws.Range("J116").Value = 1419
$ws.Range("K116").Value = 3694.0002
$ws.Range("L116").Value = 4257
$ws.Range("M116").Value = -252.0001999999999
$ws.Range("N116").Value = -11141

$ws.Range("H131").Value = 39801.31
$ws.Range("J131").Value = 1972.0834
$ws.Range("L131").Value = 5916.2502
$ws.Range("N131").Value = -15996.2502

$ws.Range("H132").Value = 2440.1765
$ws.Range("I132").Value = 1345.1428
$ws.Range("K132").Value = 12106.2852
$ws.Range("M132").Value = -9576.2852

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4216.8687
$ws.Range("I132").Value = 3347.318
$ws.Range("J132").Value = 5412.5
$ws.Range("K132").Value = 10041.954
$ws.Range("L132").Value = 16237.5
$ws.Range("M132").Value = -7511.954000000002
$ws.Range("N132").Value = -21297.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 27809.416
$ws.Range("I7").Value = 8127
$ws.Range("J7").Value = 37650.625
$ws.Range("K7").Value = 8127
$ws.Range("L7").Value = 37650.625
$ws.Range("M7").Value = -8015
$ws.Range("N7").Value = -37874.625

$ws.Range("H22").Value = 855
$ws.Range("I22").Value = 553.75
$ws.Range("J22").Value = 1156.25
$ws.Range("K22").Value = 553.75
$ws.Range("L22").Value = 1156.25
$ws.Range("M22").Value = -258.75
$ws.Range("N22").Value = -1746.25

$ws.Range("H27").Value = 855
$ws.Range("I27").Value = 553.75
$ws.Range("J27").Value = 1156.25
$ws.Range("K27").Value = 553.75
$ws.Range("L27").Value = 1156.25
$ws.Range("M27").Value = -446.75
$ws.Range("N27").Value = -1370.25

$ws.Range("H68").Value = 4054.5
$ws.Range("I68").Value = 3638.8
$ws.Range("K68").Value = 3638.8
$ws.Range("M68").Value = -2889.8

$ws.Range("H71").Value = 4054.5
$ws.Range("I71").Value = 3638.8
$ws.Range("K71").Value = 18194
$ws.Range("M71").Value = -14450

$ws.Range("H126").Value = 27809.416
$ws.Range("I126").Value = 8127
$ws.Range("J126").Value = 37650.625
$ws.Range("K126").Value = 24381
$ws.Range("L126").Value = 112951.875
$ws.Range("M126").Value = -21911
$ws.Range("N126").Value = -117891.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9907.936
$ws.Range("I62").Value = 8703.450000000001
$ws.Range("J62").Value = 12097.909
$ws.Range("K62").Value = 8703.450000000001
$ws.Range("L62").Value = 12097.909
$ws.Range("M62").Value = -8079.450000000001
$ws.Range("N62").Value = -13345.909

$ws.Range("H65").Value = 9907.936
$ws.Range("I65").Value = 8703.450000000001
$ws.Range("J65").Value = 12097.909
$ws.Range("K65").Value = 43517.25
$ws.Range("L65").Value = 60489.545
$ws.Range("M65").Value = -40397.25
$ws.Range("N65").Value = -66729.545

$ws.Range("H126").Value = 2656.2727
$ws.Range("I126").Value = 2433.625
$ws.Range("K126").Value = 7300.875
$ws.Range("M126").Value = -4830.875

$ws.Range("H132").Value = 1176255
$ws.Range("I132").Value = 875.8889
$ws.Range("K132").Value = 2627.6667
$ws.Range("M132").Value = -97.66670000000022

$ws.Range("H135").Value = 83985.664
$ws.Range("J135").Value = 83985.664
$ws.Range("L135").Value = 83985.664
$ws.Range("N135").Value = -94125.664

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 3248
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 3248
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 37332.668
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws.Range("H141").Value = 142791.6
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 142791.6
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 142791.6
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -153151.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 79099.89999999999
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 6800
$ws.Range("I17").Value = 5001
$ws.Range("J17").Value = 7399.6665
$ws.Range("K17").Value = 5001
$ws.Range("L17").Value = 7399.6665
$ws.Range("M17").Value = -4831
$ws.Range("N17").Value = -7739.6665
